$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "311.89") are stored as literal text, matching the source data,
# instead of being auto-coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.938.45'
$ws.Range("E2").Value = '  -0.32%  '

# Row 3
$ws.Range("D3").Value = '1.856.49'
$ws.Range("E3").Value = '  -1.56%  '

# Row 4
$ws.Range("E4").Value = '  +0.25%  '

# Row 5
$ws.Range("D5").Value = '311.89'
$ws.Range("E5").Value = '  -0.38%  '

# Row 6
$ws.Range("E6").Value = '  +0.14%  '

# Row 7
$ws.Range("D7").Value = '0.5086'
$ws.Range("E7").Value = '  +1.82%  '

# Row 8
$ws.Range("D8").Value = '0.3796'
$ws.Range("E8").Value = '  -1.66%  '

# Row 9
$ws.Range("D9").Value = '0.08339'
$ws.Range("E9").Value = '  -8.64%  '

# Row 10
$ws.Range("D10").Value = '1.104'
$ws.Range("E10").Value = '  -1.60%  '

# Row 11
$ws.Range("D11").Value = '41.30'
$ws.Range("E11").Value = '  -0.93%  '

# Row 12
$ws.Range("D12").Value = '6.188'
$ws.Range("E12").Value = '  -2.17%  '

# Row 13
$ws.Range("D13").Value = '1.866.18'
$ws.Range("E13").Value = '  -0.90%  '

# Row 14
$ws.Range("D14").Value = '20.37'
$ws.Range("E14").Value = '  -1.63%  '

# Row 15
$ws.Range("D15").Value = '7.162'
$ws.Range("E15").Value = '  -1.60%  '

# Row 16
$ws.Range("E16").Value = '  +0.26%  '

# Row 17
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").Value = '  -0.84%  '

# Row 18
$ws.Range("D18").Value = '90.08'
$ws.Range("E18").Value = '  -1.49%  '

# Row 19
$ws.Range("D19").Value = '0.06620'
$ws.Range("E19").Value = '  -0.14%  '

# Row 20
$ws.Range("D20").Value = '17.80'
$ws.Range("E20").Value = '  -0.20%  '

# Row 21
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.15%  '

# Row 22
$ws.Range("D22").Value = '5.998'
$ws.Range("E22").Value = '  -3.27%  '

# Row 23
$ws.Range("D23").Value = '27.981.43'
$ws.Range("E23").Value = '  -0.32%  '

# Row 24
$ws.Range("D24").Value = '11.06'
$ws.Range("E24").Value = '  -2.33%  '

# Row 25
$ws.Range("D25").Value = '2.258'
$ws.Range("E25").Value = '  -2.75%  '

# Row 26
$ws.Range("D26").Value = '2.557'
$ws.Range("E26").Value = '  +0.79%  '

# Row 27
$ws.Range("D27").Value = '2.080.67'
$ws.Range("E27").Value = '  -0.94%  '

# Row 28
$ws.Range("D28").Value = '157.10'
$ws.Range("E28").Value = '  -0.63%  '

# Row 29
$ws.Range("D29").Value = '20.41'
$ws.Range("E29").Value = '  -1.43%  '

# Row 30
$ws.Range("D30").Value = '125.77'
$ws.Range("E30").Value = '  -0.54%  '

# Row 31
$ws.Range("D31").Value = '0.1053'
$ws.Range("E31").Value = '  +0.23%  '

# Row 32
$ws.Range("D32").Value = '1.034'
$ws.Range("E32").Value = '  -3.05%  '

# Row 33
$ws.Range("D33").Value = '5.562'
$ws.Range("E33").Value = '  -0.22%  '

# Row 34
$ws.Range("D34").Value = '3.599'
$ws.Range("E34").Value = '  -0.02%  '

# Row 35
$ws.Range("D35").Value = '9.614'
$ws.Range("E35").Value = '  +2.45%  '

# Row 36
$ws.Range("D36").Value = '0.02413'
$ws.Range("E36").Value = '  +0.83%  '

# Row 37
$ws.Range("D37").Value = '0.06495'
$ws.Range("E37").Value = '  -0.70%  '

# Row 38
$ws.Range("D38").Value = '0.2149'
$ws.Range("E38").Value = '  -1.33%  '

# Row 39
$ws.Range("E39").Value = '  -0.14%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '1.231'
$ws.Range("E40").Value = '  -6.50%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6353'
$ws.Range("E41").Value = '  -0.78%  '

# Row 42
$ws.Range("D42").Value = '11.23'
$ws.Range("E42").Value = '  -2.58%  '

# Row 43
$ws.Range("D43").Value = '4.841'
$ws.Range("E43").Value = '  -1.82%  '

# Row 44
$ws.Range("D44").Value = '0.6034'
$ws.Range("E44").Value = '  +0.29%  '

# Row 45
$ws.Range("D45").Value = '12.97'
$ws.Range("E45").Value = '  -2.12%  '

# Row 46
$ws.Range("D46").Value = '1.282'
$ws.Range("E46").Value = '  -1.40%  '

# Row 47
$ws.Range("E47").Value = '  -0.40%  '

# Row 48
$ws.Range("D48").Value = '1.977'
$ws.Range("E48").Value = '  -0.56%  '

# Row 49
$ws.Range("D49").Value = '1.207'
$ws.Range("E49").Value = '  +0.59%  '

# Row 50
$ws.Range("D50").Value = '120.79'
$ws.Range("E50").Value = '  +0.10%  '

# Row 51
$ws.Range("D51").Value = '79.42'
$ws.Range("E51").Value = '  +1.06%  '

